# Updated cryptos list on Fri Jun 16 21:09:42 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'26.360.46"
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = "'  +3.16%  "
$ws.Range('E2').ClearFormats()
# Row 3
$ws.Range('D3').Value = "'1.724.07"
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = "'  +3.17%  "
$ws.Range('E3').ClearFormats()
# Row 4
$ws.Range('D4').Value = "'0.9992"
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = "'  +0.27%  "
$ws.Range('E4').ClearFormats()
# Row 5
$ws.Range('D5').Value = "'240.12"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = "'  +1.14%  "
$ws.Range('E5').ClearFormats()
# Row 6
$ws.Range('D6').Value = "'1.000"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = "'  +0.23%  "
$ws.Range('E6').ClearFormats()
# Row 7
$ws.Range('D7').Value = "'0.4722"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = "'  -1.81%  "
$ws.Range('E7').ClearFormats()
# Row 8
$ws.Range('E8').Value = "'  +0.33%  "
$ws.Range('E8').ClearFormats()
# Row 9
$ws.Range('D9').Value = "'0.06244"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = "'  +1.28%  "
$ws.Range('E9').ClearFormats()
# Row 10
$ws.Range('D10').Value = "'1.717.82"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = "'  +2.82%  "
$ws.Range('E10').ClearFormats()
# Row 11
$ws.Range('D11').Value = "'0.07085"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = "'  -0.13%  "
$ws.Range('E11').ClearFormats()
# Row 12
$ws.Range('D12').Value = "'15.37"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = "'  +3.36%  "
$ws.Range('E12').ClearFormats()
# Row 13
$ws.Range('D13').Value = "'0.5930"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = "'  -1.64%  "
$ws.Range('E13').ClearFormats()
# Row 14
$ws.Range('E14').Value = "'  -0.22%  "
$ws.Range('E14').ClearFormats()
# Row 15
$ws.Range('D15').Value = "'76.43"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = "'  +2.41%  "
$ws.Range('E15').ClearFormats()
# Row 16
$ws.Range('D16').Value = "'0.9997"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = "'  +0.16%  "
$ws.Range('E16').ClearFormats()
# Row 17
$ws.Range('D17').Value = "'1.000"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = "'  +0.34%  "
$ws.Range('E17').ClearFormats()
# Row 18
$ws.Range('D18').Value = "'26.364.07"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = "'  +3.27%  "
$ws.Range('E18').ClearFormats()
# Row 19
$ws.Range('D19').Value = "'0.000006811"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = "'  +0.03%  "
$ws.Range('E19').ClearFormats()
# Row 20
$ws.Range('D20').Value = "'11.60"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = "'  +1.10%  "
$ws.Range('E20').ClearFormats()
# Row 21
$ws.Range('D21').Value = "'1.939.28"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = "'  +3.47%  "
$ws.Range('E21').ClearFormats()
# Row 22
$ws.Range('D22').Value = "'4.565"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = "'  +1.94%  "
$ws.Range('E22').ClearFormats()
# Row 23
$ws.Range('D23').Value = "'8.773"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = "'  +0.84%  "
$ws.Range('E23').ClearFormats()
# Row 24
$ws.Range('D24').Value = "'5.341"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = "'  -0.75%  "
$ws.Range('E24').ClearFormats()
# Row 25
$ws.Range('D25').Value = "'134.91"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = "'  +0.43%  "
$ws.Range('E25').ClearFormats()
# Row 26
$ws.Range('D26').Value = "'15.28"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = "'  +1.07%  "
$ws.Range('E26').ClearFormats()
# Row 27
$ws.Range('B27').Value = "'BitcoinCash"
$ws.Range('B27').ClearFormats()
$ws.Range('C27').Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range('C27').ClearFormats()
$ws.Range('D27').Value = "'108.83"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = "'  +3.96%  "
$ws.Range('E27').ClearFormats()
# Row 28
$ws.Range('B28').Value = "'Toncoin"
$ws.Range('B28').ClearFormats()
$ws.Range('C28').Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range('C28').ClearFormats()
$ws.Range('D28').Value = "'1.411"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = "'  +0.64%  "
$ws.Range('E28').ClearFormats()
# Row 29
$ws.Range('D29').Value = "'1.768"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = "'  +3.83%  "
$ws.Range('E29').ClearFormats()
# Row 30
$ws.Range('D30').Value = "'4.027"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = "'  +1.38%  "
$ws.Range('E30').ClearFormats()
# Row 31
$ws.Range('D31').Value = "'3.703"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = "'  +0.40%  "
$ws.Range('E31').ClearFormats()
# Row 32
$ws.Range('D32').Value = "'0.07770"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = "'  +1.21%  "
$ws.Range('E32').ClearFormats()
# Row 33
$ws.Range('D33').Value = "'0.04470"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = "'  +2.40%  "
$ws.Range('E33').ClearFormats()
# Row 34
$ws.Range('E34').Value = "'  +0.06%  "
$ws.Range('E34').ClearFormats()
# Row 35
$ws.Range('D35').Value = "'0.9796"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = "'  +2.99%  "
$ws.Range('E35').ClearFormats()
# Row 36
$ws.Range('D36').Value = "'0.6224"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = "'  -0.03%  "
$ws.Range('E36').ClearFormats()
# Row 37
$ws.Range('D37').Value = "'115.70"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = "'  +17.99%  "
$ws.Range('E37').ClearFormats()
# Row 38
$ws.Range('D38').Value = "'0.9232"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = "'  +7.15%  "
$ws.Range('E38').ClearFormats()
# Row 39
$ws.Range('D39').Value = "'2.424"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = "'  -7.49%  "
$ws.Range('E39').ClearFormats()
# Row 40
$ws.Range('D40').Value = "'1.915"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = "'  +2.57%  "
$ws.Range('E40').ClearFormats()
# Row 41
$ws.Range('E41').Value = "'  +0.32%  "
$ws.Range('E41').ClearFormats()
# Row 42
$ws.Range('D42').Value = "'0.01480"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = "'  -1.66%  "
$ws.Range('E42').ClearFormats()
# Row 43
$ws.Range('D43').Value = "'5.370"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = "'  +14.81%  "
$ws.Range('E43').ClearFormats()
# Row 44
$ws.Range('D44').Value = "'0.3832"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = "'  +1.06%  "
$ws.Range('E44').ClearFormats()
# Row 45
$ws.Range('D45').Value = "'0.1169"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = "'  +4.28%  "
$ws.Range('E45').ClearFormats()
# Row 46
$ws.Range('D46').Value = "'6.286"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = "'  +0.80%  "
$ws.Range('E46').ClearFormats()
# Row 47
$ws.Range('D47').Value = "'0.05292"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = "'  +0.74%  "
$ws.Range('E47').ClearFormats()
# Row 48
$ws.Range('D48').Value = "'30.66"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = "'  +3.59%  "
$ws.Range('E48').ClearFormats()
# Row 49
$ws.Range('D49').Value = "'7.670"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = "'  +3.96%  "
$ws.Range('E49').ClearFormats()
# Row 50
$ws.Range('D50').Value = "'0.3400"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = "'  +1.33%  "
$ws.Range('E50').ClearFormats()
# Row 51
$ws.Range('D51').Value = "'1.221"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = "'  +1.77%  "
$ws.Range('E51').ClearFormats()
